$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano test statistics (DM_Stat, column C) and
# p-values (P_Value, column D) following the correction described in
# the commit message ("Correccion a Diebold Mariano").

$ws.Range("C2").Value = 0.4578811476443206
$ws.Range("D2").Value = 0.6515281364624679

$ws.Range("C3").Value = 2.101788371738125
$ws.Range("D3").Value = 0.04724539954097984

$ws.Range("C4").Value = 0.6067288523277496
$ws.Range("D4").Value = 0.5502415079792056

$ws.Range("C5").Value = 1.581505926820458
$ws.Range("D5").Value = 0.1280336269331432

$ws.Range("C6").Value = 2.362863452307251
$ws.Range("D6").Value = 0.02739224927889183

$ws.Range("C7").Value = 0.2839498404626359
$ws.Range("D7").Value = 0.7791040682120898

$ws.Range("C8").Value = 1.398589365051952
$ws.Range("D8").Value = 0.1758755565913757

$ws.Range("C9").Value = -1.490592639585562
$ws.Range("D9").Value = 0.1502677460341615

$ws.Range("C10").Value = -0.6264032473782998
$ws.Range("D10").Value = 0.5374965201911031

$ws.Range("C11").Value = 1.035510361585338
$ws.Range("D11").Value = 0.3116795417668587
